# Auto-generated edit script: updates currentAveragePrice / Leve price / profit
# columns (H:N) for the rows changed in the source commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2129.7727
$ws.Range("I96").Value = 1469.5555
$ws.Range("J96").Value = 2586.8462
$ws.Range("K96").Value = 4408.666499999999
$ws.Range("L96").Value = 7760.5386
$ws.Range("M96").Value = -3035.666499999999
$ws.Range("N96").Value = -10506.5386
$ws.Range("H125").Value = 1966.8572
$ws.Range("I125").Value = 1516
$ws.Range("J125").Value = 2147.2
$ws.Range("K125").Value = 13644
$ws.Range("L125").Value = 19324.8
$ws.Range("M125").Value = -11184
$ws.Range("N125").Value = -24244.8
$ws.Range("H129").Value = 554624.25
$ws.Range("I129").Value = 406
$ws.Range("J129").Value = 648736.75
$ws.Range("K129").Value = 1218
$ws.Range("L129").Value = 1946210.25
$ws.Range("M129").Value = 3782
$ws.Range("N129").Value = -1956210.25
$ws.Range("H137").Value = 1716.5483
$ws.Range("I137").Value = 1530.1111
$ws.Range("J137").Value = 2975
$ws.Range("K137").Value = 4590.3333
$ws.Range("L137").Value = 8925
$ws.Range("M137").Value = -2040.3333
$ws.Range("N137").Value = -14025

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1364.6364
$ws.Range("I2").Value = 735.1667
$ws.Range("J2").Value = 2120
$ws.Range("K2").Value = 735.1667
$ws.Range("L2").Value = 2120
$ws.Range("M2").Value = -622.1667
$ws.Range("N2").Value = -2346
$ws.Range("H45").Value = 1293.6471
$ws.Range("I45").Value = 1188.909
$ws.Range("J45").Value = 1485.6666
$ws.Range("K45").Value = 1188.909
$ws.Range("L45").Value = 1485.6666
$ws.Range("M45").Value = -811.9090000000001
$ws.Range("N45").Value = -2239.6666
$ws.Range("H97").Value = 493
$ws.Range("I97").Value = 491.25
$ws.Range("K97").Value = 491.25
$ws.Range("M97").Value = 4.75
$ws.Range("H110").Value = 1177.5
$ws.Range("I110").Value = 866.9375
$ws.Range("J110").Value = 1591.5834
$ws.Range("K110").Value = 866.9375
$ws.Range("L110").Value = 1591.5834
$ws.Range("M110").Value = 1178.0625
$ws.Range("N110").Value = -5681.5834
$ws.Range("H116").Value = 1364.6364
$ws.Range("I116").Value = 735.1667
$ws.Range("J116").Value = 2120
$ws.Range("K116").Value = 735.1667
$ws.Range("L116").Value = 2120
$ws.Range("M116").Value = 1558.8333
$ws.Range("N116").Value = -6708

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1364.6364
$ws.Range("I3").Value = 735.1667
$ws.Range("J3").Value = 2120
$ws.Range("K3").Value = 735.1667
$ws.Range("L3").Value = 2120
$ws.Range("M3").Value = -621.1667
$ws.Range("N3").Value = -2348
$ws.Range("H94").Value = 831.4167
$ws.Range("I94").Value = 512.125
$ws.Range("K94").Value = 512.125
$ws.Range("M94").Value = -61.125
$ws.Range("H99").Value = 1061.4286
$ws.Range("I99").Value = 919.8333
$ws.Range("K99").Value = 919.8333
$ws.Range("M99").Value = 578.1667
$ws.Range("H105").Value = 2240.7693
$ws.Range("I105").Value = 1701
$ws.Range("J105").Value = 3455.25
$ws.Range("K105").Value = 1701
$ws.Range("L105").Value = 3455.25
$ws.Range("M105").Value = 46
$ws.Range("N105").Value = -6949.25
$ws.Range("H107").Value = 1035.4736
$ws.Range("I107").Value = 1072.125
$ws.Range("J107").Value = 840
$ws.Range("K107").Value = 1072.125
$ws.Range("L107").Value = 840
$ws.Range("M107").Value = 847.875
$ws.Range("N107").Value = -4680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6063409
$ws.Range("I31").Value = 2950
$ws.Range("J31").Value = 66668000
$ws.Range("K31").Value = 2950
$ws.Range("L31").Value = 66668000
$ws.Range("M31").Value = -2655
$ws.Range("N31").Value = -66668590
$ws.Range("H34").Value = 6063409
$ws.Range("I34").Value = 2950
$ws.Range("J34").Value = 66668000
$ws.Range("K34").Value = 2950
$ws.Range("L34").Value = 66668000
$ws.Range("M34").Value = -2748
$ws.Range("N34").Value = -66668404
$ws.Range("H99").Value = 2392.25
$ws.Range("I99").Value = 1882.7778
$ws.Range("J99").Value = 2809.0908
$ws.Range("K99").Value = 1882.7778
$ws.Range("L99").Value = 2809.0908
$ws.Range("M99").Value = -384.7778000000001
$ws.Range("N99").Value = -5805.0908
$ws.Range("H107").Value = 501.03226
$ws.Range("I107").Value = 437.08334
$ws.Range("K107").Value = 437.08334
$ws.Range("M107").Value = 1482.91666
$ws.Range("H126").Value = 2392.25
$ws.Range("I126").Value = 1882.7778
$ws.Range("J126").Value = 2809.0908
$ws.Range("K126").Value = 5648.3334
$ws.Range("L126").Value = 8427.2724
$ws.Range("M126").Value = -3178.3334
$ws.Range("N126").Value = -13367.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 26750
$ws.Range("I116").Value = 34000
$ws.Range("J116").Value = 22400
$ws.Range("K116").Value = 102000
$ws.Range("L116").Value = 67200
$ws.Range("M116").Value = -98558
$ws.Range("N116").Value = -74084
$ws.Range("H131").Value = 1142228.2
$ws.Range("J131").Value = 1848986
$ws.Range("L131").Value = 5546958
$ws.Range("N131").Value = -5557038

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8043.5713
$ws.Range("I80").Value = 2326.25
$ws.Range("J80").Value = 15666.667
$ws.Range("K80").Value = 2326.25
$ws.Range("L80").Value = 15666.667
$ws.Range("M80").Value = -1328.25
$ws.Range("N80").Value = -17662.667
$ws.Range("H83").Value = 8043.5713
$ws.Range("I83").Value = 2326.25
$ws.Range("J83").Value = 15666.667
$ws.Range("K83").Value = 11631.25
$ws.Range("L83").Value = 78333.33499999999
$ws.Range("M83").Value = -6639.25
$ws.Range("N83").Value = -88317.33499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1845.2778
$ws.Range("J46").Value = 1857.125
$ws.Range("L46").Value = 1857.125
$ws.Range("N46").Value = -2233.125
$ws.Range("H55").Value = 352
$ws.Range("I55").Value = 320
$ws.Range("J55").Value = 426.66666
$ws.Range("K55").Value = 320
$ws.Range("L55").Value = 426.66666
$ws.Range("M55").Value = -147
$ws.Range("N55").Value = -772.66666
$ws.Range("H61").Value = 30304186
$ws.Range("I61").Value = 1314.2858
$ws.Range("J61").Value = 83334210
$ws.Range("K61").Value = 1314.2858
$ws.Range("L61").Value = 83334210
$ws.Range("M61").Value = -1112.2858
$ws.Range("N61").Value = -83334614
$ws.Range("H113").Value = 30304186
$ws.Range("I113").Value = 1314.2858
$ws.Range("J113").Value = 83334210
$ws.Range("K113").Value = 1314.2858
$ws.Range("L113").Value = 83334210
$ws.Range("M113").Value = 855.7141999999999
$ws.Range("N113").Value = -83338550
$ws.Range("H132").Value = 4535.5
$ws.Range("I132").Value = 4833.3335
$ws.Range("J132").Value = 4312.125
$ws.Range("K132").Value = 14500.0005
$ws.Range("L132").Value = 12936.375
$ws.Range("M132").Value = -11970.0005
$ws.Range("N132").Value = -17996.375
$ws.Range("H136").Value = 3813.2927
$ws.Range("I136").Value = 4330.484
$ws.Range("J136").Value = 2210
$ws.Range("K136").Value = 12991.452
$ws.Range("L136").Value = 6630
$ws.Range("M136").Value = -10441.452
$ws.Range("N136").Value = -11730

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2413.6155
$ws.Range("I62").Value = 2056
$ws.Range("J62").Value = 2637.125
$ws.Range("K62").Value = 2056
$ws.Range("L62").Value = 2637.125
$ws.Range("M62").Value = -1432
$ws.Range("N62").Value = -3885.125
$ws.Range("H65").Value = 2413.6155
$ws.Range("I65").Value = 2056
$ws.Range("J65").Value = 2637.125
$ws.Range("K65").Value = 10280
$ws.Range("L65").Value = 13185.625
$ws.Range("M65").Value = -7160
$ws.Range("N65").Value = -19425.625
$ws.Range("H81").Value = 2278
$ws.Range("I81").Value = 2383
$ws.Range("J81").Value = 1700.5
$ws.Range("K81").Value = 4766
$ws.Range("L81").Value = 3401
$ws.Range("M81").Value = -3705
$ws.Range("N81").Value = -5523
$ws.Range("H84").Value = 2278
$ws.Range("I84").Value = 2383
$ws.Range("J84").Value = 1700.5
$ws.Range("K84").Value = 23830
$ws.Range("L84").Value = 17005
$ws.Range("M84").Value = -18526
$ws.Range("N84").Value = -27613
$ws.Range("H96").Value = 41667970
$ws.Range("I96").Value = 41667970
$ws.Range("K96").Value = 41667970
$ws.Range("M96").Value = -41666597
$ws.Range("H100").Value = 1014.2857
$ws.Range("I100").Value = 1300
$ws.Range("J100").Value = 900
$ws.Range("K100").Value = 2600
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -2059
$ws.Range("N100").Value = -2882
$ws.Range("H107").Value = 218
$ws.Range("I107").Value = 216
$ws.Range("J107").Value = 225
$ws.Range("K107").Value = 648
$ws.Range("L107").Value = 675
$ws.Range("M107").Value = 1272
$ws.Range("N107").Value = -4515
$ws.Range("H136").Value = 4440.222
$ws.Range("I136").Value = 4667.125
$ws.Range("J136").Value = 2625
$ws.Range("K136").Value = 2625
$ws.Range("L136").Value = 7875
$ws.Range("M136").Value = -11451.375
$ws.Range("N136").Value = -12975
